# Fixed update to excel issue
#
# 1. Rename the "Requested quantity" header on the "Weekly Quantity" and
#    "Monthly Trend" sheets.
# 2. Add a new "PO Forecast" sheet (after "Monthly Trend") with forecast data.

$wb = $excel.ActiveWorkbook

# --- 1. Rename headers -----------------------------------------------------
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- 2. Add the "PO Forecast" sheet ----------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add($null, $lastSheet)
$wsForecast.Name = "PO Forecast"

# Header row
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# Copy the existing bold/bordered/centered header format from the
# "Weekly Quantity" sheet so the new header cells reuse the same style
# instead of minting a near-duplicate one.
$wsWeekly.Range("A1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)

# Data rows: ds (weekly date), PO_Forecast, yhat_lower, yhat_upper
$data = @(
    @(44934.99999999999, 5, 4.001163381691227, 6.099961382413466),
    @(44962.99999999999, 4, 2.459234343851328, 4.531606930367685),
    @(44997.99999999999, 2, 0.5540808934768763, 2.616564313341809),
    @(45004.99999999999, 1, 0.03086031738956808, 2.169716736066599),
    @(45011.99999999999, 1, -0.4153471377008456, 1.742947990279307),
    @(45018.99999999999, 0, -0.7576046847258536, 1.304049511378408),
    @(45025.99999999999, 0, -1.102985529252547, 0.8690116667525281),
    @(45032.99999999999, 0, -1.408487595766245, 0.4935300731179382),
    @(45039.99999999999, 0, -1.912375494341977, 0.1629894999411809),
    @(45046.99999999999, 0, -2.260603076318876, -0.2900810958469332),
    @(45053.99999999999, 0, -2.63778691857335, -0.650806872451043),
    @(45060.99999999999, 0, -3.241644426282959, -0.9535677506254686),
    @(45067.99999999999, 0, -3.512612571783672, -1.402434332705747)
)

$row = 2
foreach ($r in $data) {
    $wsForecast.Cells.Item($row, 1).Value = $r[0]
    $wsForecast.Cells.Item($row, 2).Value = $r[1]
    $wsForecast.Cells.Item($row, 3).Value = $r[2]
    $wsForecast.Cells.Item($row, 4).Value = $r[3]
    $row = $row + 1
}

# Reuse the existing date-formatted style (style 2) for the "ds" column,
# same way the header row reused style 1 above.
$wsWeekly.Range("A2").Copy()
$wsForecast.Range("A2:A14").PasteSpecial(-4122)

# Restore the original active sheet/selection so this edit doesn't shift
# the workbook's active-tab state as a side effect.
$wsWeekly.Activate() | Out-Null
$wsWeekly.Range("A1").Select() | Out-Null
